$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text (add units / rename Notes column) ---
$ws.Range("A1").Value = "h (km)"
$ws.Range("B1").Value = "rhoMin (g/km3)"
$ws.Range("C1").Value = "rhoMax (g/km3)"
$ws.Range("D1").Value = "Notes"

# --- Insert a new data row right below the header for h=0 (Wikipedia reference row) ---
[void]$ws.Rows(2).Insert()
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1225000000000
$ws.Range("C2").Value = 1225000000000
$ws.Range("B2:C2").NumberFormat = "0.00E+00"
$ws.Range("D2").Value = "Wikipedia"

# --- Fix the old mistyped "93..08" text value (now shifted to row 17) with the real number ---
$ws.Range("C17").Value = 93.08

# --- Append a final row representing density at (effectively) infinite altitude ---
$ws.Range("A53").Value = 10000000000
$ws.Range("A53").NumberFormat = "0.00E+00"
$ws.Range("B53").Value = 0
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = "Inf altitude"

# --- Column widths for the newly visible / widened columns ---
$ws.Columns("B").ColumnWidth = 14.1640625
$ws.Columns("C").ColumnWidth = 14.5
$ws.Columns("D").ColumnWidth = 12.5

# --- Update selection to match final saved state ---
[void]$ws.Range("D53").Select()

Write-Output "edit complete"
